$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D updates to remain text (values look numeric, e.g. "591.97")
# by briefly tagging the cell as Text format, then restoring the default
# "Normal" style afterwards so the cell keeps no explicit style index,
# matching the original inline-string cells.

$ws.Range("D2").Value = "64.164.84"
$ws.Range("E2").Value = "  -2.17%  "

$ws.Range("D3").Value = "3.111.57"
$ws.Range("E3").Value = "  -3.15%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.51%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "

$ws.Range("D9").Value = "3.110.07"
$ws.Range("E9").Value = "  -3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.74%  "

$ws.Range("E12").Value = "  -4.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.98%  "

$ws.Range("E14").Value = "  -6.01%  "

$ws.Range("D15").Value = "3.626.05"
$ws.Range("E15").Value = "  -3.15%  "

$ws.Range("E16").Value = "  -1.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.66%  "

$ws.Range("D18").Value = "64.099.11"
$ws.Range("E18").Value = "  -1.66%  "

$ws.Range("D19").Value = "3.112.97"
$ws.Range("E19").Value = "  -3.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.61%  "

$ws.Range("E22").Value = "  -8.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.05%  "

$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("E30").Value = "  -4.20%  "

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  -4.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.74%  "

$ws.Range("E34").Value = "  -4.85%  "

$ws.Range("D35").Value = "0.0₃0840"
$ws.Range("E35").Value = "  -6.75%  "

$ws.Range("E36").Value = "  -2.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("E41").Value = "  -2.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "440.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.02%  "

$ws.Range("E43").Value = "  -3.89%  "

$ws.Range("E44").Value = "  -5.56%  "

$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "

$ws.Range("D47").Value = "2.839.71"
$ws.Range("E47").Value = "  -3.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("E51").Value = "  -3.99%  "

